# Weekly refresh of the Femacal de La Calera - Espinaca series:
# a new week's record is inserted at the top of the data block (row 521),
# and every existing record shifts down one row (the previously-last
# record, old row 592, becomes row 593).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 521, pushing rows 521:592 down
# to 522:593 (dimension grows from R592 to R593 automatically).
$ws.Rows("521:521").Insert()

# Populate the newly inserted row 521 with the new week's data.
$ws.Range("A521").Value = 3
$ws.Range("B521").Value = "Femacal de La Calera"
$ws.Range("C521").Value = "Coquimbo"
$ws.Range("D521").Value = 45154
$ws.Range("E521").Value = 5
$ws.Range("F521").Value = 100112012
$ws.Range("G521").Value = "Espinaca"
$ws.Range("H521").Value = "Sin especificar"
$ws.Range("I521").Value = "Primera"
$ws.Range("J521").Value = 80
$ws.Range("K521").Value = 4500
$ws.Range("L521").Value = 4500
$ws.Range("M521").Value = 4500
$ws.Range("N521").Value = '$/docena de atados (3 kilos)'
$ws.Range("O521").Value = "Provincia de Quillota"
$ws.Range("P521").Value = 1500
$ws.Range("Q521").Value = 3
$ws.Range("R521").Value = "Hortaliza"
